# Replace the table style applied to every table in the deck from the
# custom "Table_0" style ({D9C2FC74-120F-4CD5-8DDD-B74F4DED4CBD}) to the
# built-in "No Style, No Grid" table style
# ({59AF0FED-7032-46D9-A783-0CF860946D73}), matching what a user does by
# picking that style from the Table Design > Table Styles gallery.
#
# Table.Style is read-only in this host (assigning to it raises "Table
# styles cannot be assigned through a property - call
# Table.ApplyStyle(\"{GUID}\") instead"), so Table.ApplyStyle(...) is used.

$p = $ppt.ActivePresentation

$oldStyleId = "{D9C2FC74-120F-4CD5-8DDD-B74F4DED4CBD}"
$newStyleId = "{59AF0FED-7032-46D9-A783-0CF860946D73}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $table = $shape.Table
            if ($table.Style -eq $oldStyleId) {
                $table.ApplyStyle($newStyleId)
            }
        }
    }
}
